# excel : utilisation de la bdd administration pour récup date début / fin
#
# Update the "period" label in H3 (end date of the training period) and
# the now-unused placeholder value in Q7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H3: "du 09/09/2019 au 08/09/2020" -> "du 09/09/2019 au 27/05/2020"
$ws.Range("H3").Value = "du 09/09/2019 au 27/05/2020"

# Q7: "21" -> "0"
$ws.Range("Q7").Value = "0"
